# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handoff XLIFF files were (re)generated for the 4 "Ready for
# handoff" / low-priority files, bumping their priority to "ht" and
# refreshing the relevant handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 4 affected files
# (rows 4-7) moves from 12:39:13 to 12:39:37.
$overview.Range("G4").Value = "2016-09-05 12:39:37"
$overview.Range("G5").Value = "2016-09-05 12:39:37"
$overview.Range("G6").Value = "2016-09-05 12:39:37"
$overview.Range("G7").Value = "2016-09-05 12:39:37"

# zh-cn sheet: Priority goes from "low" to "ht" (handed off), and the
# Latest Handoff Datetime is refreshed from 12:39:00 to 12:39:31.
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("E6").Value = "ht"
$zhcn.Range("E7").Value = "ht"

$zhcn.Range("H4").Value = "2016-09-05 12:39:31"
$zhcn.Range("H5").Value = "2016-09-05 12:39:31"
$zhcn.Range("H6").Value = "2016-09-05 12:39:31"
$zhcn.Range("H7").Value = "2016-09-05 12:39:31"

# de-de sheet: Priority goes from "low" to "ht" (handed off), and the
# Latest Handoff Datetime is refreshed from 12:39:13 to 12:39:37 (shares
# the same generation timestamp as the Overview sheet).
$dede.Range("E4").Value = "ht"
$dede.Range("E5").Value = "ht"
$dede.Range("E6").Value = "ht"
$dede.Range("E7").Value = "ht"

$dede.Range("H4").Value = "2016-09-05 12:39:37"
$dede.Range("H5").Value = "2016-09-05 12:39:37"
$dede.Range("H6").Value = "2016-09-05 12:39:37"
$dede.Range("H7").Value = "2016-09-05 12:39:37"
